# The deck had an extra "ML Modeling Framework / Model Evaluation" slide
# (slide 5) that only contained a short bullet list (no charts/images) -
# essentially a duplicate/placeholder of the following "Model Evaluation"
# slides which already carry the full bullet list plus the confusion
# matrix / metrics screenshots. Remove that redundant slide; the rest of
# the deck (and the slides after it) keep their order.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$s.Delete()
